# Rename header cells from "constants.kN" -> "parameters.kN"
# (commit message: "constants -> parameters")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("conditions")

$ws.Range("B1").Value = "parameters.k1"
$ws.Range("C1").Value = "parameters.k2"
$ws.Range("D1").Value = "parameters.k3"
